# "Generate Report for Archive" -- refresh the localization-status report:
# every cell still showing the old "Ready for handoff" status is now
# "In Translation" (Overview!E2:F4 plus the Status column on each
# per-locale sheet), and the Status-ish columns are re-sized to fit the
# shorter text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newWidth = 13.4101845877511

# --- Overview sheet: "zh-cn" (col E) and "de-de" (col F) status columns ---
$overview = $wb.Worksheets.Item("Overview")
$overviewRange = $overview.Range("E2:F4")
foreach ($cell in $overviewRange.Cells) {
    if ($cell.Value() -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}
$overview.Columns("E").ColumnWidth = $newWidth
$overview.Columns("F").ColumnWidth = $newWidth

# --- Per-locale sheets: "Status" column (col C) ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $statusRange = $ws.Range("C2:C4")
    foreach ($cell in $statusRange.Cells) {
        if ($cell.Value() -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
    $ws.Columns("C").ColumnWidth = $newWidth
}
